$d = $word.ActiveDocument

# --- Title / TitleChar: drop the condensed-tracking / kerning tweak ---
# (Font dialog "Character Spacing: Normal", kerning off -> Spacing/Kerning = 0)
$titlePara = $d.Styles("Title")
$titlePara.Font.Spacing = 0
$titlePara.Font.Kerning = 0

$titleChar = $d.Styles("TitleChar")
$titleChar.Font.Spacing = 0
$titleChar.Font.Kerning = 0

# --- Author: now based on Title (picks up the centering from there),
#     with its own run size pinned back down to 12pt ---
$author = $d.Styles("Author")
$author.BaseStyle = $d.Styles("Title")
$author.Font.Size = 12
$author.Font.SizeBi = 12

# --- Date: same treatment ---
$date = $d.Styles("Date")
$date.BaseStyle = $d.Styles("Title")
$date.Font.Size = 12
$date.Font.SizeBi = 12
